$d = $word.ActiveDocument

# The document ends with a "REFERENCES USED :" heading paragraph followed
# by three more paragraphs (two hyperlinks and a "youtube.com ..." line).
# The last of those paragraphs also carries the (hidden) _GoBack bookmark.
# The edit removes all of the reference text, leaving a single empty
# paragraph - the former "REFERENCES USED :" paragraph - which keeps its
# paragraph mark/formatting and inherits the _GoBack bookmark.

$refParaIndex = 0
$lastParaIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "REFERENCES USED :") {
        $refParaIndex = $i
    }
    if ($t -eq "youtube.com – various videos.") {
        $lastParaIndex = $i
    }
}

if ($refParaIndex -gt 0 -and $lastParaIndex -ge $refParaIndex) {
    $refPara = $d.Paragraphs($refParaIndex)

    # 1) Remove the "REFERENCES USED :" run text, keeping the paragraph
    #    itself (and its paragraph mark / pPr formatting) intact.
    $clearEnd = $refPara.Range.End - 1
    if ($clearEnd -gt $refPara.Range.Start) {
        $d.Range($refPara.Range.Start, $clearEnd).Delete()
    }

    if ($lastParaIndex -gt $refParaIndex) {
        $firstLinkPara = $d.Paragraphs($refParaIndex + 1)
        $lastRefPara = $d.Paragraphs($lastParaIndex)

        # 2) Delete the hyperlink paragraph(s) and the trailing
        #    "youtube.com ..." paragraph entirely, including their
        #    paragraph marks, which merges them out of the document and
        #    leaves the (now empty) REFERENCES USED paragraph as the last
        #    paragraph of the body.
        $d.Range($firstLinkPara.Range.Start, $lastRefPara.Range.End).Delete()
    }

    # 3) The _GoBack bookmark used to live in the final (now deleted)
    #    paragraph; recreate it, zero-length, at the end of the paragraph
    #    that survives so it keeps marking "last edit position".
    $survivor = $d.Paragraphs($refParaIndex)
    $d.Bookmarks.Add("_GoBack", $survivor.Range)
}
